# Insert a new data row at row 163 (pushes existing rows 163..238 down to
# 164..239, which is exactly what the target diff shows: every row from the
# old row 163 onward keeps its original data but moves one row lower, the
# former last row (238) becomes row 239, and the dimension grows from
# A1:R238 to A1:R239).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(163).EntireRow.Insert()

# Populate the freshly inserted row 163 with its new record.
$ws.Range("A163").Value = 3
$ws.Range("B163").Value = "Femacal de La Calera"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 44523
$ws.Range("E163").Value = 5
$ws.Range("F163").Value = 100112043
$ws.Range("G163").Value = "Pepino ensalada"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 140
$ws.Range("K163").Value = 6500
$ws.Range("L163").Value = 7000
$ws.Range("M163").Value = 6714
$ws.Range("N163").Value = "$/caja 70 unidades"
$ws.Range("O163").Value = "Región de Arica y Parinacota"
$ws.Range("P163").Value = 96
$ws.Range("Q163").Value = 70
$ws.Range("R163").Value = "Hortaliza"
